$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Protokol mit den Faellen des Servers ergaenzt:
# zwei neue Protokoll-Zeilen (Loginversuch / Register - Username Passwort)
$ws.Range("B10").Value = "Loginversuch - Username Passwort"
$ws.Range("B11").Value = "Register - Username Passwort"

# Move/update the active selection to B12
$ws.Range("B12").Select()

# Resize the workbook window (best-effort; matches author's window resize)
$win = $excel.ActiveWindow
$win.Width = 38400
$win.Height = 12225
